# Add a new worksheet named "PARSE" at the end of the workbook, populate it
# with a small parsed-data table, and make it the active sheet/tab - mirrors
# introducing a "parsedSheet" type that is read once and cached.

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "PARSE"

# Header row
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
$ws.Range("D1").Value = "d"

# Data rows
$ws.Range("A2").Value = "abc"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1.1111

$ws.Range("A3").Value = "def"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1.233

$ws.Range("A4").Value = "jkl"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3.24342

# Date-valued column, formatted as a short date (reuses the workbook's
# existing built-in date style) and filled with a volatile TODAY() formula.
$ws.Range("D2:D4").NumberFormat = "m/d/yy"
$ws.Range("D2").Formula = "=TODAY()"
$ws.Range("D3:D4").Formula = "=TODAY()"

$ws.Columns.Item(4).AutoFit()

# Leave the new sheet selected/active, matching the position a user would be
# in right after typing the data in (D2:D4 highlighted, D4 last-entered).
[void]$ws.Activate()
[void]$ws.Range("D2:D4").Select()
